$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.142.30'
$ws.Range("E2").Value = '  -2.57%  '
$ws.Range("D3").Value = '1.719.11'
$ws.Range("E3").Value = '  -2.82%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.17%  '
$ws.Range("D5").Value = '309.21'
$ws.Range("E5").Value = '  -5.70%  '
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("D7").Value = '0.4730'
$ws.Range("E7").Value = '  +5.47%  '
$ws.Range("D8").Value = '0.3432'
$ws.Range("E8").Value = '  -3.62%  '
$ws.Range("D9").Value = '42.18'
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").Value = '0.07276'
$ws.Range("E10").Value = '  -2.19%  '
$ws.Range("D11").Value = '1.045'
$ws.Range("E11").Value = '  -4.76%  '
$ws.Range("D12").Value = '0.9998'
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("D13").Value = '19.92'
$ws.Range("E13").Value = '  -4.93%  '
$ws.Range("D14").Value = '5.880'
$ws.Range("E14").Value = '  -2.40%  '
$ws.Range("D15").Value = '1.716.47'
$ws.Range("E15").Value = '  -3.59%  '
$ws.Range("D16").Value = '6.888'
$ws.Range("E16").Value = '  -4.99%  '
$ws.Range("D17").Value = '89.33'
$ws.Range("E17").Value = '  -4.19%  '
$ws.Range("D18").Value = '0.00001041'
$ws.Range("E18").Value = '  -1.76%  '
$ws.Range("D19").Value = '0.06358'
$ws.Range("E19").Value = '  -1.24%  '
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").Value = '16.54'
$ws.Range("E21").Value = '  -3.43%  '
$ws.Range("D22").Value = '5.625'
$ws.Range("E22").Value = '  -2.68%  '
$ws.Range("D23").Value = '27.164.28'
$ws.Range("E23").Value = '  -2.68%  '
$ws.Range("D24").Value = '10.88'
$ws.Range("E24").Value = '  -3.61%  '
$ws.Range("E25").Value = '  +0.57%  '
$ws.Range("D26").Value = '156.07'
$ws.Range("E26").Value = '  -4.28%  '
$ws.Range("D27").Value = '19.52'
$ws.Range("E27").Value = '  -4.10%  '
$ws.Range("D28").Value = '1.920.53'
$ws.Range("E28").Value = '  -3.08%  '
$ws.Range("D29").Value = '2.092'
$ws.Range("E29").Value = '  -2.77%  '
$ws.Range("D30").Value = '119.61'
$ws.Range("E30").Value = '  -4.24%  '
$ws.Range("D31").Value = '1.019'
$ws.Range("E31").Value = '  -7.98%  '
$ws.Range("D32").Value = '0.09156'
$ws.Range("E32").Value = '  -0.32%  '
$ws.Range("D33").Value = '3.596'
$ws.Range("D34").Value = '5.328'
$ws.Range("E34").Value = '  -5.11%  '
$ws.Range("D35").Value = '0.02208'
$ws.Range("E35").Value = '  -3.66%  '
$ws.Range("D36").Value = '0.05822'
$ws.Range("E36").Value = '  -4.33%  '
$ws.Range("D37").Value = '11.03'
$ws.Range("E37").Value = '  -7.05%  '
$ws.Range("D38").Value = '0.1998'
$ws.Range("E38").Value = '  -4.96%  '
$ws.Range("D39").Value = '4.745'
$ws.Range("E39").Value = '  -4.30%  '
$ws.Range("D40").Value = '1.396'
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("D41").Value = '0.5904'
$ws.Range("E41").Value = '  -6.68%  '
$ws.Range("D42").Value = '1.118'
$ws.Range("E42").Value = '  -5.61%  '
$ws.Range("D43").Value = '7.488'
$ws.Range("E43").Value = '  -5.09%  '
$ws.Range("D44").Value = '12.52'
$ws.Range("E44").Value = '  -6.43%  '
$ws.Range("D45").Value = '0.5665'
$ws.Range("E45").Value = '  -4.15%  '
$ws.Range("D46").Value = '3.568'
$ws.Range("E46").Value = '  -4.67%  '
$ws.Range("D47").Value = '117.56'
$ws.Range("E47").Value = '  -3.91%  '
$ws.Range("D48").Value = '1.846'
$ws.Range("E48").Value = '  -5.67%  '
$ws.Range("D49").Value = '0.06653'
$ws.Range("E49").Value = '  -3.64%  '
$ws.Range("E50").Value = '  -4.32%  '
$ws.Range("D51").Value = '1.000'
$ws.Range("E51").Value = '  +0.03%  '
